$wb = $excel.ActiveWorkbook

# --- Add "booking" sheet (sheetId=2) right after "data1" ---
$ws1 = $wb.Worksheets.Item("data1")
$booking = $wb.Worksheets.Add($null, $ws1)
$booking.Name = "booking"

$booking.Range("A1").Value = "username"
$booking.Range("B1").Value = "password"
$booking.Range("A2").Value = "admin"
$booking.Range("B2").Value = "password123"

# --- Add "booking_data" sheet (sheetId=3) right after "booking" ---
$bookingData = $wb.Worksheets.Add($null, $booking)
$bookingData.Name = "booking_data"

$bookingData.Range("A1").Value = "firstname"
$bookingData.Range("B1").Value = "lastname"
$bookingData.Range("C1").Value = "totalprice"
$bookingData.Range("D1").Value = "depositpaid"
$bookingData.Range("E1").Value = "checkin"
$bookingData.Range("F1").Value = "checkout"
$bookingData.Range("G1").Value = "additionalneeds"

$bookingData.Range("A2").Value = "Akram"
$bookingData.Range("B2").Value = "Wasim"
$bookingData.Range("C2").Value = 2000

# "false" must land as literal text, not a boolean -- use a text formula assignment
$bookingData.Range("D2").Formula = "'false"

$bookingData.Range("E2").Value = [datetime]"2023-01-01"
$bookingData.Range("E2").NumberFormat = "mm-dd-yy"
$bookingData.Range("F2").Value = [datetime]"2023-01-02"
$bookingData.Range("F2").NumberFormat = "mm-dd-yy"

$bookingData.Range("G2").Value = "JungleTour-"

# --- Column sizing (approximate autofit to content) ---
$booking.Range("A1:B2").Columns.AutoFit()
$bookingData.Range("A1:G2").Columns.AutoFit()

# --- Selections matching the edited workbook state ---
$ws1.Range("G12").Select()
$booking.Range("A1:B2").Select()

$bookingData.Range("G2").Select()
$bookingData.Activate()
